$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new reading was recorded for 2026/01/13 (09:00) which belongs between the
# existing 2026/01/13 06:00 row (625) and the 2026/12/29 13:00 row (626).
# Copy the row immediately above (625: date/weekday text already correct,
# "2026/01/13" / "火") and use Copy+Insert so the whole block from row 626
# downward shifts down by one row while the new row inherits the exact same
# (unstyled) cell formatting as the rest of the table. Only the hour value
# (column C) then needs to be corrected to 9.
$ws.Range("A625:D625").Copy()
$ws.Range("A626:D626").Insert()
$ws.Range("C626").Value = 9

Write-Host "Inserted row 626 (2026/01/13, 火, 9, 201); rows 626-667 shifted to 627-668."
